# "Generate Report for Handback" -- refresh the handback status report:
#  - the en-US sync status flips to "not in sync" (shown on the Overview
#    sheet's per-language columns and on each language sheet's Status column)
#  - the ad1ecfcc file's Correspond Handback DateTime is refreshed for both
#    the zh-cn and de-de languages
#  - the report's generated columns are widened a bit to fit the new data

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Handed back: not in sync with en-US"

# Overview sheet: zh-cn / de-de status columns for both rows
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Per-language detail sheets: Status column (C) for both rows
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Refresh the "Correspond Handback DateTime" for the ad1ecfcc handback file
$zhcn.Range("K3").Value = "2016-11-29 05:14:50"
$dede.Range("K3").Value = "2016-11-29 05:15:09"

# Widen the generated-content columns to fit the refreshed report
$overview.Range("E1:F1").ColumnWidth = 32.6666666666667
$overview.Range("F1").ColumnWidth = 32.6666666666667
$zhcn.Range("C1").ColumnWidth = 32.6666666666667
$dede.Range("C1").ColumnWidth = 32.6666666666667
